# Update regression coefficient values across the four sheets
# (exponential, weibull, lognormal, loglogistic) to reflect the fix
# described in the commit: the function was using the wrong data
# structure/format, so the computed coefficients changed.

$wb = $excel.ActiveWorkbook

# --- Sheet: exponential ---
$ws1 = $wb.Worksheets.Item("exponential")
$ws1.Range("C2").Value = 0.0606465275779165
$ws1.Range("D2").Value = 0.0527586331715868
$ws1.Range("E2").Value = 0.0540637390173801

# --- Sheet: weibull ---
$ws2 = $wb.Worksheets.Item("weibull")
$ws2.Range("B3").Value = 0.0589862330212438
$ws2.Range("B4").Value = -0.0113584882966721
$ws2.Range("C4").Value = 0.0452961847902155

$ws2.Range("B8").Value = 0.0619592603432408
$ws2.Range("B9").Value = -0.0128179425420927
$ws2.Range("C9").Value = 0.0406511128175405

$ws2.Range("B13").Value = 0.0405087071958447
$ws2.Range("B14").Value = 0.00275245948638132
$ws2.Range("C14").Value = 0.0525507587379983

# --- Sheet: lognormal ---
$ws3 = $wb.Worksheets.Item("lognormal")
$ws3.Range("B3").Value = 0.0617677799350474
$ws3.Range("B4").Value = -0.0019259391554576
$ws3.Range("C4").Value = 0.0445890580851386

$ws3.Range("B8").Value = 0.0562954952310476
$ws3.Range("B9").Value = -0.00309717690119501
$ws3.Range("C9").Value = 0.042669435266441

$ws3.Range("B13").Value = 0.0495164494110894
$ws3.Range("B14").Value = 0.00854800860377524
$ws3.Range("C14").Value = 0.0460967621500772

# --- Sheet: loglogistic ---
$ws4 = $wb.Worksheets.Item("loglogistic")
$ws4.Range("B3").Value = 0.0640912155632968
$ws4.Range("B4").Value = -0.000982529676603334
$ws4.Range("C4").Value = 0.0487655237850791

$ws4.Range("B8").Value = 0.0567696803186672
$ws4.Range("B9").Value = 0.00141827357831441
$ws4.Range("C9").Value = 0.0480355167606601

$ws4.Range("B13").Value = 0.0458233048611451
$ws4.Range("B14").Value = 0.00436994056411851
$ws4.Range("C14").Value = 0.052409913253889
